$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.527.85'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.823.08'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.10'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5105'
$ws.Range("E7").Value = '  -5.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3947'
$ws.Range("E8").Value = '  -2.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08174'
$ws.Range("E9").Value = '  +6.49%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.109'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.65'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.344'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.11'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.514'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '1.821.46'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001131'
$ws.Range("E17").Value = '  +3.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.48'
$ws.Range("E18").Value = '  +3.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06663'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.83'
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.0000'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.095'
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("D23").Value = '28.546.12'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.32'
$ws.Range("E26").Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.19'
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").Value = '2.029.74'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.403'
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.92'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.765'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07065'
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02354'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.246'
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.842'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.180'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.399'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.50'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5925'
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.08'
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.185'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06893'
$ws.Range("E51").Value = '  +0.30%  '
